$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1232.75
$ws.Range("I2").Value = 1045.1
$ws.Range("J2").Value = 1545.5
$ws.Range("K2").Value = 1045.1
$ws.Range("L2").Value = 1545.5
$ws.Range("M2").Value = -932.0999999999999
$ws.Range("N2").Value = -1771.5
$ws.Range("H17").Value = 753
$ws.Range("J17").Value = 755.5476
$ws.Range("L17").Value = 2266.6428
$ws.Range("N17").Value = -2602.6428
$ws.Range("H38").Value = 25.75
$ws.Range("I38").Value = 25.75
$ws.Range("K38").Value = 77.25
$ws.Range("M38").Value = 294.75
$ws.Range("H86").Value = 1973.5
$ws.Range("I86").Value = 1664
$ws.Range("J86").Value = 2437.75
$ws.Range("K86").Value = 1664
$ws.Range("L86").Value = 2437.75
$ws.Range("M86").Value = -541
$ws.Range("N86").Value = -4683.75
$ws.Range("H89").Value = 1973.5
$ws.Range("I89").Value = 1664
$ws.Range("J89").Value = 2437.75
$ws.Range("K89").Value = 8320
$ws.Range("L89").Value = 12188.75
$ws.Range("M89").Value = -2704
$ws.Range("N89").Value = -23420.75
$ws.Range("H96").Value = 8929584
$ws.Range("I96").Value = 23809828
$ws.Range("J96").Value = 1437.8
$ws.Range("K96").Value = 71429484
$ws.Range("L96").Value = 4313.4
$ws.Range("M96").Value = -71428111
$ws.Range("N96").Value = -7059.4
$ws.Range("H103").Value = 681.41174
$ws.Range("I103").Value = 410
$ws.Range("K103").Value = 1230
$ws.Range("M103").Value = -644
$ws.Range("H106").Value = 3336.5557
$ws.Range("I106").Value = 3336.5557
$ws.Range("K106").Value = 3336.5557
$ws.Range("M106").Value = -2705.5557
$ws.Range("H116").Value = 17260.043
$ws.Range("I116").Value = 17308.643
$ws.Range("K116").Value = 17308.643
$ws.Range("M116").Value = -13866.643
$ws.Range("H125").Value = 76981.75
$ws.Range("I125").Value = 300032
$ws.Range("J125").Value = 2631.6667
$ws.Range("K125").Value = 2700288
$ws.Range("L125").Value = 23685.0003
$ws.Range("M125").Value = -2697828
$ws.Range("N125").Value = -28605.0003
$ws.Range("H131").Value = 3680.1765
$ws.Range("I131").Value = 3298.7273
$ws.Range("J131").Value = 4379.5
$ws.Range("K131").Value = 9896.1819
$ws.Range("L131").Value = 13138.5
$ws.Range("M131").Value = -4856.1819
$ws.Range("N131").Value = -23218.5
$ws.Range("H141").Value = 1800.5
$ws.Range("I141").Value = 1649.75
$ws.Range("J141").Value = 2102
$ws.Range("K141").Value = 4949.25
$ws.Range("L141").Value = 6306
$ws.Range("M141").Value = 230.75
$ws.Range("N141").Value = -16666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").Value = ""
$ws.Range("H61").Value = 4355.6055
$ws.Range("I61").Value = 1039.7858
$ws.Range("K61").Value = 1039.7858
$ws.Range("M61").Value = -827.7858000000001
$ws.Range("H122").Value = 1436.9584
$ws.Range("I122").Value = 1465.2106
$ws.Range("J122").Value = 1329.6
$ws.Range("K122").Value = 4395.6318
$ws.Range("L122").Value = 3988.8
$ws.Range("M122").Value = -1945.6318
$ws.Range("N122").Value = -8888.8
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
$ws.Range("H132").Value = 1346.0149
$ws.Range("I132").Value = 801.90247
$ws.Range("K132").Value = 2405.70741
$ws.Range("M132").Value = 124.29259
$ws.Range("H136").Value = 4355.6055
$ws.Range("I136").Value = 1039.7858
$ws.Range("K136").Value = 3119.3574
$ws.Range("M136").Value = -569.3574000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1965.1818
$ws.Range("I94").Value = 1487.7142
$ws.Range("K94").Value = 1487.7142
$ws.Range("M94").Value = -1036.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7695376.5
$ws.Range("I31").Value = 12500987
$ws.Range("J31").Value = 6400.4
$ws.Range("K31").Value = 12500987
$ws.Range("L31").Value = 6400.4
$ws.Range("M31").Value = -12500692
$ws.Range("N31").Value = -6990.4
$ws.Range("H34").Value = 7695376.5
$ws.Range("I34").Value = 12500987
$ws.Range("J34").Value = 6400.4
$ws.Range("K34").Value = 12500987
$ws.Range("L34").Value = 6400.4
$ws.Range("M34").Value = -12500785
$ws.Range("N34").Value = -6804.4
$ws.Range("H132").Value = 168819
$ws.Range("I132").Value = 250752.25
$ws.Range("K132").Value = 752256.75
$ws.Range("M132").Value = -749726.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2173.8572
$ws.Range("I3").Value = 2173.8572
$ws.Range("K3").Value = 6521.571599999999
$ws.Range("M3").Value = -6409.571599999999
$ws.Range("H133").Value = 2802.7
$ws.Range("I133").Value = 2802.7
$ws.Range("K133").Value = 8408.099999999999
$ws.Range("M133").Value = -3348.099999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8762.04
$ws.Range("I70").Value = 9503.75
$ws.Range("J70").Value = 7443.4443
$ws.Range("K70").Value = 9503.75
$ws.Range("L70").Value = 7443.4443
$ws.Range("M70").Value = -9233.75
$ws.Range("N70").Value = -7983.4443
$ws.Range("H73").Value = 8762.04
$ws.Range("I73").Value = 9503.75
$ws.Range("J73").Value = 7443.4443
$ws.Range("K73").Value = 9503.75
$ws.Range("L73").Value = 7443.4443
$ws.Range("M73").Value = -8567.75
$ws.Range("N73").Value = -9315.4443
$ws.Range("H102").Value = 57241.727
$ws.Range("I102").Value = 62466.5
$ws.Range("K102").Value = 62466.5
$ws.Range("M102").Value = -60844.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1854.75
$ws.Range("J22").Value = 1750.2858
$ws.Range("L22").Value = 1750.2858
$ws.Range("N22").Value = -2340.2858
$ws.Range("H27").Value = 1854.75
$ws.Range("J27").Value = 1750.2858
$ws.Range("L27").Value = 1750.2858
$ws.Range("N27").Value = -1964.2858
$ws.Range("H46").Value = 7227.3335
$ws.Range("J46").Value = 7227.3335
$ws.Range("L46").Value = 7227.3335
$ws.Range("N46").Value = -7603.3335
$ws.Range("H68").Value = 4754.3
$ws.Range("I68").Value = 3499.5
$ws.Range("J68").Value = 5068
$ws.Range("K68").Value = 3499.5
$ws.Range("L68").Value = 5068
$ws.Range("M68").Value = -2750.5
$ws.Range("N68").Value = -6566
$ws.Range("H71").Value = 4754.3
$ws.Range("I71").Value = 3499.5
$ws.Range("J71").Value = 5068
$ws.Range("K71").Value = 17497.5
$ws.Range("L71").Value = 25340
$ws.Range("M71").Value = -13753.5
$ws.Range("N71").Value = -32828
$ws.Range("H82").Value = 1279.4706
$ws.Range("J82").Value = 1370
$ws.Range("L82").Value = 1370
$ws.Range("N82").Value = -2092
$ws.Range("H85").Value = 1279.4706
$ws.Range("J85").Value = 1370
$ws.Range("L85").Value = 1370
$ws.Range("N85").Value = -3866
$ws.Range("H93").Value = 1617
$ws.Range("I93").Value = 1305.1818
$ws.Range("K93").Value = 1305.1818
$ws.Range("M93").Value = -57.18180000000007
$ws.Range("H132").Value = 1597.875
$ws.Range("I132").Value = 789.86664
$ws.Range("J132").Value = 2944.5557
$ws.Range("K132").Value = 2369.59992
$ws.Range("L132").Value = 8833.667099999999
$ws.Range("M132").Value = 160.4000800000003
$ws.Range("N132").Value = -13893.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1201.7142
$ws.Range("I100").Value = 918.5
$ws.Range("K100").Value = 1837
$ws.Range("M100").Value = -1296
$ws.Range("H122").Value = 130108.055
$ws.Range("I122").Value = 193233
$ws.Range("K122").Value = 579699
$ws.Range("M122").Value = -577249
